$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "All" with "Combined" in column A for the rows that contain it
foreach ($r in 2,5,8,11,14,17) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "All") {
        $cell.Value = "Combined"
    }
}
